$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "dct:creator"
$ws.Cells.Item(13, 2).Value = "https://orcid.org/0000-0002-0454-4289"
$ws.Cells.Item(13, 3).Value = "Hannah Mihai"
